# Add new columns I ("I0") and J ("IF") to the sheet, mirroring the style
# of the existing header row and filling in the per-row data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns, copying the style of the existing
# header cell (H1) so the new header cells match formatting.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for rows 2..42 (columns I and J)
$data = @(
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(6,7),
    @(9,9),
    @(7,7),
    @(7,8),
    @(9,9),
    @(7,7),
    @(7,8),
    @(7,8),
    @(8,8),
    @(8,8),
    @(1,2),
    @(4,5),
    @(7,8),
    @(8,9),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(5,6),
    @(6,7),
    @(7,7),
    @(5,6),
    @(6,7),
    @(7,7),
    @(4,5),
    @(9,9)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

$wb.Save()
